# Auto-generated Excel COM-interop script applying numeric updates
# described by the authoritative OOXML diff (market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1247.3158
$ws.Range("I5").Value = 182.66667
$ws.Range("K5").Value = 182.66667
$ws.Range("M5").Value = -67.66667000000001
$ws.Range("H8").Value = 10000694
$ws.Range("I8").Value = 14285777
$ws.Range("K8").Value = 42857331
$ws.Range("M8").Value = -42857192
$ws.Range("H32").Value = 30001730
$ws.Range("I32").Value = 40001680
$ws.Range("K32").Value = 40001680
$ws.Range("M32").Value = -40001354
$ws.Range("H40").Value = 5725.5454
$ws.Range("I40").Value = 6091.1665
$ws.Range("J40").Value = 5286.8
$ws.Range("K40").Value = 6091.1665
$ws.Range("L40").Value = 5286.8
$ws.Range("M40").Value = -5916.1665
$ws.Range("N40").Value = -5636.8
$ws.Range("H43").Value = 2232.6667
$ws.Range("I43").Value = 1891
$ws.Range("K43").Value = 1891
$ws.Range("M43").Value = -1822
$ws.Range("H62").Value = 4118.294
$ws.Range("I62").Value = 3836.2727
$ws.Range("K62").Value = 3836.2727
$ws.Range("M62").Value = -3212.2727
$ws.Range("H65").Value = 4118.294
$ws.Range("I65").Value = 3836.2727
$ws.Range("K65").Value = 19181.3635
$ws.Range("M65").Value = -16061.3635
$ws.Range("H92").Value = 291
$ws.Range("I92").Value = 113.76471
$ws.Range("J92").Value = 1797.5
$ws.Range("K92").Value = 113.76471
$ws.Range("L92").Value = 1797.5
$ws.Range("M92").Value = 1134.23529
$ws.Range("N92").Value = -4293.5
$ws.Range("H112").Value = 2759.1191
$ws.Range("J112").Value = 2910.8057
$ws.Range("L112").Value = 8732.417099999999
$ws.Range("N112").Value = -10948.4171
$ws.Range("H116").Value = 12110.667
$ws.Range("I116").Value = 15296
$ws.Range("J116").Value = 8925.333000000001
$ws.Range("K116").Value = 15296
$ws.Range("L116").Value = 8925.333000000001
$ws.Range("M116").Value = -11854
$ws.Range("N116").Value = -15809.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13215031
$ws.Range("I2").Value = 13910480
$ws.Range("K2").Value = 13910480
$ws.Range("M2").Value = -13910367
$ws.Range("H39").Value = 4000
$ws.Range("I39").Value = 4000
$ws.Range("K39").Value = 4000
$ws.Range("M39").Value = -3480
$ws.Range("H40").Value = 42498.5
$ws.Range("I40").Value = 49999
$ws.Range("J40").Value = 34998
$ws.Range("K40").Value = 49999
$ws.Range("L40").Value = 34998
$ws.Range("M40").Value = -49823
$ws.Range("N40").Value = -35350
$ws.Range("H42").Value = 30983.166
$ws.Range("J42").Value = 30983.166
$ws.Range("L42").Value = 30983.166
$ws.Range("M42").Value = -31955.166
$ws.Range("H116").Value = 13215031
$ws.Range("I116").Value = 13910480
$ws.Range("K116").Value = 13910480
$ws.Range("M116").Value = -13908186
$ws.Range("H132").Value = 2167.13
$ws.Range("I132").Value = 2035.1794
$ws.Range("J132").Value = 2634.9546
$ws.Range("K132").Value = 6105.5382
$ws.Range("L132").Value = 7904.8638
$ws.Range("M132").Value = -3575.5382
$ws.Range("N132").Value = -12964.8638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13215031
$ws.Range("I3").Value = 13910480
$ws.Range("K3").Value = 13910480
$ws.Range("M3").Value = -13910366
$ws.Range("H20").Value = 3591.3142
$ws.Range("I20").Value = 4234.9165
$ws.Range("J20").Value = 2187.0908
$ws.Range("K20").Value = 4234.9165
$ws.Range("L20").Value = 2187.0908
$ws.Range("M20").Value = -3987.9165
$ws.Range("N20").Value = -2681.0908
$ws.Range("H64").Value = 13334959
$ws.Range("I64").Value = 20834976
$ws.Range("K64").Value = 20834976
$ws.Range("M64").Value = -20834751
$ws.Range("H67").Value = 13334959
$ws.Range("I67").Value = 20834976
$ws.Range("K67").Value = 20834976
$ws.Range("M67").Value = -20834196
$ws.Range("H94").Value = 958.6111
$ws.Range("I94").Value = 569.619
$ws.Range("J94").Value = 1503.2
$ws.Range("K94").Value = 569.619
$ws.Range("L94").Value = 1503.2
$ws.Range("M94").Value = -118.619
$ws.Range("N94").Value = -2405.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 19765.166
$ws.Range("I22").Value = 19765.166
$ws.Range("K22").Value = 19765.166
$ws.Range("M22").Value = -19415.166
$ws.Range("H99").Value = 9088.27
$ws.Range("I99").Value = 6939.24
$ws.Range("K99").Value = 6939.24
$ws.Range("M99").Value = -5441.24
$ws.Range("H125").Value = 40314.4
$ws.Range("J125").Value = 40314.4
$ws.Range("L125").Value = 40314.4
$ws.Range("N125").Value = -45234.4
$ws.Range("H126").Value = 9088.27
$ws.Range("I126").Value = 6939.24
$ws.Range("K126").Value = 20817.72
$ws.Range("M126").Value = -18347.72
$ws.Range("H132").Value = 6231.1963
$ws.Range("I132").Value = 1418.3125
$ws.Range("J132").Value = 14337.105
$ws.Range("K132").Value = 4254.9375
$ws.Range("L132").Value = 43011.315
$ws.Range("M132").Value = -1724.9375
$ws.Range("N132").Value = -48071.315
$ws.Range("H134").Value = 2885.5593
$ws.Range("I134").Value = 2412.0557
$ws.Range("J134").Value = 7999.4
$ws.Range("K134").Value = 7236.1671
$ws.Range("L134").Value = 23998.2
$ws.Range("M134").Value = -4701.1671
$ws.Range("N134").Value = -29068.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 225
$ws.Range("I18").Value = 225
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 675
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -506
$ws.Range("N18").ClearContents()
$ws.Range("H34").Value = 397.42856
$ws.Range("I34").Value = 62.1
$ws.Range("J34").Value = 702.2727
$ws.Range("K34").Value = 186.3
$ws.Range("L34").Value = 2106.8181
$ws.Range("M34").Value = -102.3
$ws.Range("N34").Value = -2274.8181
$ws.Range("I46").Value = 2960.8
$ws.Range("J46").Value = 324306340
$ws.Range("K46").Value = 8882.400000000001
$ws.Range("L46").Value = 972919020
$ws.Range("M46").Value = -8791.400000000001
$ws.Range("N46").Value = -972919202
$ws.Range("H100").Value = 5466.6665
$ws.Range("J100").Value = 6900
$ws.Range("L100").Value = 20700
$ws.Range("N100").Value = -22322
$ws.Range("H108").Value = 597
$ws.Range("I108").Value = 597
$ws.Range("K108").Value = 1791
$ws.Range("M108").Value = 1089
$ws.Range("H110").Value = 8537.200000000001
$ws.Range("I110").Value = 8949.75
$ws.Range("K110").Value = 26849.25
$ws.Range("M110").Value = -22759.25
$ws.Range("H119").Value = 2555.5
$ws.Range("J119").Value = 9025
$ws.Range("L119").Value = 27075
$ws.Range("N119").Value = -36751
$ws.Range("H123").Value = 1300
$ws.Range("I123").Value = 1300
$ws.Range("K123").Value = 3900
$ws.Range("M123").Value = -1450
$ws.Range("H124").Value = 6268.636
$ws.Range("I124").Value = 3994.5
$ws.Range("K124").Value = 11983.5
$ws.Range("M124").Value = -7073.5
$ws.Range("H125").Value = 2136
$ws.Range("I125").Value = 2136
$ws.Range("K125").Value = 6408
$ws.Range("M125").Value = -1488
$ws.Range("H126").Value = 801
$ws.Range("I126").Value = 801
$ws.Range("K126").Value = 2403
$ws.Range("M126").Value = 2537
$ws.Range("H130").Value = 1100
$ws.Range("I130").Value = 1100
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 3300
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 1720
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 1839292.8
$ws.Range("I131").Value = 2674817
$ws.Range("J131").Value = 1139.2
$ws.Range("K131").Value = 8024451
$ws.Range("L131").Value = 3417.6
$ws.Range("M131").Value = -8019411
$ws.Range("N131").Value = -13497.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 8055.5
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H102").Value = 3803.8494
$ws.Range("I102").Value = 2840.0173
$ws.Range("J102").Value = 7530.6665
$ws.Range("K102").Value = 2840.0173
$ws.Range("L102").Value = 7530.6665
$ws.Range("M102").Value = -1218.0173
$ws.Range("N102").Value = -10774.6665
$ws.Range("H113").Value = 4281.1
$ws.Range("I113").Value = 1784.75
$ws.Range("J113").Value = 5945.3335
$ws.Range("K113").Value = 1784.75
$ws.Range("L113").Value = 5945.3335
$ws.Range("M113").Value = 385.25
$ws.Range("N113").Value = -10285.3335
$ws.Range("H122").Value = 835291.3
$ws.Range("I122").Value = 1430272.4
$ws.Range("K122").Value = 4290817.199999999
$ws.Range("M122").Value = -4288367.199999999
$ws.Range("H132").Value = 2736.724
$ws.Range("I132").Value = 1554.7778
$ws.Range("J132").Value = 4670.8184
$ws.Range("K132").Value = 4664.3334
$ws.Range("L132").Value = 14012.4552
$ws.Range("M132").Value = -2134.3334
$ws.Range("N132").Value = -19072.4552
$ws.Range("H136").Value = 35556.24
$ws.Range("J136").Value = 35556.24
$ws.Range("L136").Value = 106668.72
$ws.Range("N136").Value = -111768.72

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3472822
$ws.Range("I16").Value = 4310794
$ws.Range("K16").Value = 4310794
$ws.Range("M16").Value = -4310624
$ws.Range("H132").Value = 6704.091
$ws.Range("I132").Value = 5448.6
$ws.Range("J132").Value = 8356.053
$ws.Range("K132").Value = 16345.8
$ws.Range("L132").Value = 25068.159
$ws.Range("M132").Value = -13815.8
$ws.Range("N132").Value = -30128.159

